$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.428.08"
$ws.Range("E2").Value = "  +3.30%  "
$ws.Range("D3").Value = "2.648.67"
$ws.Range("E3").Value = "  +2.14%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "606.78"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.72%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "156.76"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.16%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.595"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.62%  "
$ws.Range("E9").Value = "  +8.56%  "
$ws.Range("E10").Value = "  +4.67%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.83"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.14%  "
$ws.Range("E12").Value = "  +1.56%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "29.11"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +6.12%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000189"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +21.04%  "
$ws.Range("D15").Value = "3.126.45"
$ws.Range("E15").Value = "  +2.22%  "
$ws.Range("D16").Value = "65.311.03"
$ws.Range("E16").Value = "  +3.31%  "
$ws.Range("D17").Value = "2.635.61"
$ws.Range("E17").Value = "  +1.52%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.66"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.75%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.84"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.55%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "355.51"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.87%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.29"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +5.96%  "
$ws.Range("E22").Value = "  +0.33%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "68.23"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.96%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.74"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.06%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.55"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.41%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.67"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.22%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.17"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.39%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.165"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.76%  "
$ws.Range("D29").Value = "0.0₃0960"
$ws.Range("E29").Value = "  +13.87%  "
$ws.Range("E30").Value = "  -0.53%  "
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.13"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.96%  "
$ws.Range("B32").Value = "Bittensor"
$ws.Range("C32").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "525.72"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.01%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.80"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.18%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.71"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +10.40%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.42"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.51%  "
$ws.Range("E36").Value = "  +4.19%  "
$ws.Range("E37").Value = "  +6.99%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "165.69"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "20.25"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.52%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.00"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.02%  "
$ws.Range("E41").Value = "  -0.04%  "
$ws.Range("B42").Value = "OKB"
$ws.Range("C42").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "42.41"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +7.01%  "
$ws.Range("B43").Value = "Aave"
$ws.Range("C43").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "167.07"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.01%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.11"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.30%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0610"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.99%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "23.36"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.27%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.22"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +9.72%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.654"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.87%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0253"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.70%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0991"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.96%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "19.67"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.14%  "
